# "Generate Report for Handoff" -- regenerate the localization-status report:
#   * Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#   * The zh-cn handoff timestamp advances (14:56:57 -> 14:57:50)
#   * The HO-xliff-generate timestamp advances (14:57:04 -> 14:57:55), which
#     shows up on the Overview sheet and on the de-de detail sheet
#   * The (now shorter) Status column narrows to fit its new contents on all
#     three sheets

$wb = $excel.ActiveWorkbook

$statusOld = "Handed back: in sync with en-US"
$statusNew = "Ready for handoff"

# ---- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("G2").Value = "2016-08-27 14:57:55"
$wsOverview.Columns.Item(5).ColumnWidth = 16.25
$wsOverview.Columns.Item(6).ColumnWidth = 16.25

# ---- zh-cn detail sheet ----------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("H2").Value = "2016-08-27 14:57:50"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.25

# ---- de-de detail sheet ------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("H2").Value = "2016-08-27 14:57:55"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.25
